# Delete rows 13 and 14 (even_MAG-GUT49384.fa and even_MAG-GUT6224.fa),
# which shifts the subsequent rows up by two and shrinks the used range
# from A1:Y32 down to A1:Y30.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(13).Resize(2).Delete()
